$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row (13) with the new test mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Range("D13").Value = "Productinformatie"
$logs.Range("E13").Value = "Beste klant,`nBedankt voor uw e-mail. Om uw vraag over de leverbaarheid te beantwoorden, hebben we meer informatie nodig over het specifieke product waar u naar informeert. Kunt u ons de naam of het artikelnummer van het product geven?`nZodra we deze informatie hebben, zullen we direct voor u nakijken of het product nog leverbaar is. Alvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Bedrijfsnaam] e-mailassistent"
$logs.Range("F13").Value = "2025-08-01 23:51:50"
$logs.Range("G13").Value = "Ja"
$logs.Range("H13").Value = "Nee"
$logs.Range("I13").Value = "Ja"
$logs.Range("J13").Value = "Nee"

# The multi-line text in E13 triggers an automatic row-height bump; restore
# the row to the sheet's standard (default) height, like the other rows.
$logs.Rows.Item(13).AutoFit()

# Extend the conditional-formatting ranges so they keep covering the new row
$logs.Range("D2:D13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D13"))
$logs.Range("G2:G13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G13"))
$logs.Range("H2:H13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H13"))
$logs.Range("I2:I13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I13"))
$logs.Range("J2:J13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J13"))

# --- Sheet "Dashboard": update the category summary counts ---
# Row 3 becomes "Productinformatie" with count 3
# Row 4 becomes "Intern verzoek / Actie voor medewerker" with count 2
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 3
$dash.Range("A4").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B4").Value = 2
